# Update the division problems in the table to the new values from the
# commit. Cells are addressed positionally (row, column) rather than by
# text search because several expressions (e.g. "94÷7=", "54÷4=") occur
# more than once in the document, and a text-based Find/Replace would not
# be able to distinguish between them reliably.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="96÷8="},
    @{Row=1;  Col=2; Text="90÷7="},
    @{Row=1;  Col=3; Text="81÷8="},
    @{Row=1;  Col=4; Text="54÷4="},
    @{Row=1;  Col=5; Text="54÷7="},

    @{Row=5;  Col=1; Text="32÷9="},
    @{Row=5;  Col=2; Text="38÷8="},
    @{Row=5;  Col=3; Text="88÷6="},
    @{Row=5;  Col=4; Text="26÷6="},
    @{Row=5;  Col=5; Text="40÷5="},

    @{Row=9;  Col=1; Text="18÷9="},
    @{Row=9;  Col=2; Text="26÷2="},
    @{Row=9;  Col=3; Text="17÷8="},
    @{Row=9;  Col=4; Text="74÷5="},
    @{Row=9;  Col=5; Text="60÷4="},

    @{Row=13; Col=1; Text="57÷7="},
    @{Row=13; Col=2; Text="91÷9="},
    @{Row=13; Col=3; Text="26÷9="},
    @{Row=13; Col=4; Text="52÷9="},
    @{Row=13; Col=5; Text="41÷5="},

    @{Row=17; Col=1; Text="72÷5="},
    @{Row=17; Col=2; Text="60÷3="},
    @{Row=17; Col=3; Text="91÷8="},
    @{Row=17; Col=4; Text="26÷8="},
    @{Row=17; Col=5; Text="17÷7="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $u.Text
}
